$d = $word.ActiveDocument

# Touch the run's Bold property so the engine rewrites <w:b w:val="true"/>
# as the canonical boolean-true form (keeping it bold==true, matching the
# "true" -> "on" ST_OnOff synonym change in the source diff). Use a Range
# that stops before the paragraph mark so the paragraph's own rPr is left
# untouched.
$storyLen = $d.Content.End
$r = $d.Range(0, $storyLen - 1)
$r.Font.Bold = $true

# Bump the stack-trace line numbers (Apache POI 4.1.0 -> 5.2.3 upgrade
# shifted a few JDK internal line references).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "ThreadPoolExecutor.runWorker(ThreadPoolExecutor.java:1130)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ThreadPoolExecutor.runWorker(ThreadPoolExecutor.java:1136)", 2
)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute(
    "ThreadPoolExecutor`$Worker.run(ThreadPoolExecutor.java:630)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ThreadPoolExecutor`$Worker.run(ThreadPoolExecutor.java:635)", 2
)

$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Replacement.ClearFormatting()
$find3.Execute(
    "java.lang.Thread.run(Thread.java:832)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "java.lang.Thread.run(Thread.java:833)", 2
)
